$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update U2 / U3 text values (shared strings merge: blink -> blink+templerun, templerun -> blink+sudoku)
$ws.Range("U2").Value = "blink+templerun"
$ws.Range("U3").Value = "blink+sudoku"

# Update row 2 numeric results
$ws.Range("P2").Value = 0.01
$ws.Range("Q2").Value = -1.021077917979345
$ws.Range("R2").Value = 11.67664403205391
$ws.Range("S2").Value = -8.677719827163163

# Update row 3 numeric results
$ws.Range("Q3").Value = -1.195184938981114
$ws.Range("R3").Value = 12.84974359013185
$ws.Range("S3").Value = -8.965394171833633

# Remove rows 4 through 6 (shrinking dimension to A1:U3)
$ws.Range("A4:U6").Delete()
